# Update the "Fixed" date placeholder text shown on the slide master,
# every slide layout, and the notes master from 3/4/2024 to 3/5/2024.
#
# These placeholders hold a <a:fld type="datetimeFigureOut"> whose cached
# <a:t> text PowerPoint keeps in sync with the Header/Footer date text.
# ppPlaceholderDate == 16.

$ppPlaceholderDate = 16
$oldDate = "3/4/2024"
$newDate = "3/5/2024"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $isDatePlaceholder = $false
            if ($shp.Type -eq 14) {
                if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                    $isDatePlaceholder = $true
                }
            }
            if ($isDatePlaceholder) {
                if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                    $shp.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

$p = $ppt.ActivePresentation

# 1. Slide master.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# 2. Every slide layout ("custom layout") hanging off the master.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# 3. Notes master. Its date placeholder shape can't be written directly
# through Shapes(...).TextFrame.TextRange in this host (silently
# ignored), so the write goes through HeadersFooters.DateAndTime
# instead; that object's own .Text getter is unreliable, so the
# shape's TextFrame is used to read/verify the current value.
$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    $shp = $notesMaster.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.Type -eq 14) {
            if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                    $notesMaster.HeadersFooters.DateAndTime.Text = $newDate
                }
            }
        }
    }
}
